$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column C (rows 2-14) with the same scores already present in column B
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 5

# Comments added in column D
$ws.Range("D6").Value = "**attention Radioactive Goods a le mauvais titre"
$ws.Range("D15").Value = "** radioactive goods a le mauvais titre de page, les fields dans la form qui sont nécessaires devrais être noté (field required)"

# Update the selected cell to match the saved view
$ws.Range("C5").Select()

$wb.Save()
